$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting the
# existing "Late" / "Outstanding" / "heading" / "Disbursement" columns
# one position to the right (N->O, O->P, P->Q).
$ws.Columns("N:N").Insert()

# The newly inserted column picks up the width of the column immediately
# to its left ("In Advance", column M).
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Reflect the cell Excel leaves selected after performing this edit.
$ws.Range("R9").Select()
